# Commit: "add file with mapping of variable names to product attribute descriptions"
#
# The workbook already contains a mapping table on Sheet2 (category / variable
# name / description / gsub) used to document each modelling attribute. This
# edit reworks the free-text descriptions in column D into a consistent,
# self-explanatory phrasing ("<Label>: Indicator variable, equaling 1 if
# ..., 0 otherwise ...") and highlights the updated columns (C:E) in yellow
# so readers can see which parts of the table were touched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# --- Rewrite the attribute descriptions in column D (rows 2-43) ------------
$ws.Cells.Item(2, 4).Value2  = "Screensize: Indicator variable, equaling 1 if screensize <= 24 inch, 0 otherwise (baseline: unknown screensizes)"
$ws.Cells.Item(3, 4).Value2  = "Screensize: Indicator variable, equaling 1 if screensize > 24 inch, 0 otherwise (baseline: unknown screensizes)"
$ws.Cells.Item(4, 4).Value2  = "Megapixels"
$ws.Cells.Item(5, 4).Value2  = "Digital zoom: Indicator variable, equaling 1 if available, 0 otherwise"
$ws.Cells.Item(6, 4).Value2  = "Megapixels"
$ws.Cells.Item(7, 4).Value2  = "Blueray: Indicator variable, equaling 1 if available, 0 otherwise"
$ws.Cells.Item(8, 4).Value2  = "Recording function: Indicator variable, equaling 1 if available, 0 otherwise"
$ws.Cells.Item(9, 4).Value2  = "RAM (in MB)"
$ws.Cells.Item(10, 4).Value2 = "CPU (in Mhz)"
$ws.Cells.Item(11, 4).Value2 = "Hard disk size"
$ws.Cells.Item(12, 4).Value2 = "without LEDs: Indicator variable, equaling 1 for LCD TVs without LEDs, 0 otherwise"
$ws.Cells.Item(13, 4).Value2 = "Screensize: Indicator variable, equaling 1 if screensize < 40 inch, 0 otherwise (baseline: unknown screensizes)"
$ws.Cells.Item(14, 4).Value2 = "Screensize: Indicator variable, equaling 1 if screensize >= 40 inch, 0 otherwise (baseline: unknown screensizes)"
$ws.Cells.Item(15, 4).Value2 = "3D: Indicator variable, equaling 1 if available, 0 otherwise"
$ws.Cells.Item(16, 4).Value2 = "RAM (in MB)"
$ws.Cells.Item(17, 4).Value2 = "CPU (in Mhz)"
$ws.Cells.Item(18, 4).Value2 = "Weight (in kg)"
$ws.Cells.Item(19, 4).Value2 = "Touchscreen: Indicator variable, equaling 1 if available, 0 otherwise"
$ws.Cells.Item(20, 4).Value2 = "Webcam: Indicator variable, equaling 1 if available, 0 otherwise"
$ws.Cells.Item(21, 4).Value2 = "Capacity (in liters)"
$ws.Cells.Item(22, 4).Value2 = "Functionality: Indicator variable, equaling 1 for only microwave functionality, 0 otherwise (e.g., if microwave comes with a grill)"
$ws.Cells.Item(23, 4).Value2 = "Power (in watt)"
$ws.Cells.Item(24, 4).Value2 = "Digital time controller: Indicator variable, equaling 1 if available, 0 otherwise"
$ws.Cells.Item(25, 4).Value2 = "Touchscreen: Indicator variable, equaling 1 if available, 0 otherwise"
$ws.Cells.Item(26, 4).Value2 = "Wifi: Indicator variable, equaling 1 if available, 0 otherwise"
$ws.Cells.Item(27, 4).Value2 = "Bluetooth: Indicator variable, equaling 1 if available, 0 otherwise"
$ws.Cells.Item(28, 4).Value2 = "Screensize: Indicator variable, equaling 1 if screensize < 40 inch, 0 otherwise (baseline: unknown screensizes)"
$ws.Cells.Item(29, 4).Value2 = "Screensize: Indicator variable, equaling 1 if screensize >= 40 inch, 0 otherwise (baseline: unknown screensizes)"
$ws.Cells.Item(30, 4).Value2 = "3D: Indicator variable, equaling 1 if available, 0 otherwise"
$ws.Cells.Item(31, 4).Value2 = "Freezer: Indicator variable, equaling 1 if available, 0 otherwise"
$ws.Cells.Item(32, 4).Value2 = "Number of doors"
$ws.Cells.Item(33, 4).Value2 = "Touchscreen: Indicator variable, equaling 1 if available, 0 otherwise"
$ws.Cells.Item(34, 4).Value2 = "Screensize (in inch)"
$ws.Cells.Item(35, 4).Value2 = "Wifi: Indicator variable, equaling 1 if available, 0 otherwise"
$ws.Cells.Item(36, 4).Value2 = "Memory (in GB)"
$ws.Cells.Item(37, 4).Value2 = "RAM (in MB)"
$ws.Cells.Item(38, 4).Value2 = "Weight (in kg)"
$ws.Cells.Item(39, 4).Value2 = "Screensize (in inch)"
$ws.Cells.Item(40, 4).Value2 = "Capacity (in liters)"
$ws.Cells.Item(41, 4).Value2 = "Functionality: Indicator variable, equaling 1 if washing machine comes with extra functionality; 0 if washing only"
$ws.Cells.Item(42, 4).Value2 = "Front loader: Indicator variable, equaling 1 if front loader; 0 otherwise (baseline is others)"
$ws.Cells.Item(43, 4).Value2 = "Top loader: Indicator variable, equaling 1 if top loader; 0 otherwise (baseline is others)"

# --- Highlight the variable-name / description / gsub columns in yellow ----
$ws.Range("C2:E43").Interior.Color = 65535

# --- Resize column D to fit the longer descriptions -------------------------
$ws.Columns.Item(4).AutoFit()

# --- Restore the active selection shown when the workbook is reopened ------
$ws.Range("D16").Select()
